# Update NATMI ligand-receptor stats with new TPM-derived values.
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
#          E Ligand-expressing cells, F Ligand detection rate,
#          G Ligand average expression value, H Ligand total expression value,
#          I/J Ligand derived specificity (avg/total),
#          K Receptor-expressing cells, L Receptor detection rate,
#          M Receptor average expression value, N Receptor total expression value,
#          O/P Receptor derived specificity (avg/total),
#          Q/R Edge expression weight (avg/total),
#          S/T Edge expression derived specificity (avg/total)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E2" = 3;                    "F2" = 1;                    "G2" = 0.292879
    "H2" = 0.878637;             "I2" = 0.007455471808396097; "J2" = 0.007455471808396097
    "M2" = 0.1809866666666667;   "N2" = 0.54296;               "O2" = 0.03987407676082905
    "P2" = 0.03987407676082905;  "Q2" = 0.05300719394666666;  "R2" = 0.47706474552
    "S2" = 0.0002972800551761829; "T2" = 0.0002972800551761829

    "E3" = 3;                    "F3" = 1;                    "G3" = 0.292879
    "H3" = 0.878637;             "I3" = 0.007455471808396097; "J3" = 0.007455471808396097
    "O3" = 0.1057193993302571;   "P3" = 0.1057193993302571;   "Q3" = 0.1405396477976667
    "R3" = 1.264856830179;       "S3" = 0.0007881880013073011; "T3" = 0.0007881880013073011

    "E4" = 3;                    "F4" = 1;                    "G4" = 0.292879
    "H4" = 0.878637;             "I4" = 0.007455471808396097; "J4" = 0.007455471808396097
    "M4" = 3.878113333333333;    "N4" = 11.63434;              "O4" = 0.8544065239089139
    "P4" = 0.8544065239089139;   "Q4" = 1.135817954953333;    "R4" = 10.22236159458
    "S4" = 0.006370003751912613; "T4" = 0.006370003751912613

    "I5" = 0.6729737392616156;   "J5" = 0.6729737392616155
    "M5" = 0.1809866666666667;   "N5" = 0.54296;               "O5" = 0.03987407676082905
    "P5" = 0.03987407676082905;  "Q5" = 4.784734009439999;    "R5" = 43.06260608496
    "S5" = 0.02683420653733982;  "T5" = 0.02683420653733981

    "I6" = 0.6729737392616156;   "J6" = 0.6729737392616155
    "O6" = 0.1057193993302571;   "P6" = 0.1057193993302571
    "S6" = 0.07114637947977506;  "T6" = 0.07114637947977505

    "I7" = 0.6729737392616156;   "J7" = 0.6729737392616155
    "M7" = 3.878113333333333;    "N7" = 11.63434;              "O7" = 0.8544065239089139
    "P7" = 0.8544065239089139;   "Q7" = 102.52545726276;       "R7" = 922.7291153648399
    "S7" = 0.5749931532445008;   "T7" = 0.5749931532445007

    "G8" = 12.55394366666667;    "H8" = 37.661831;             "I8" = 0.3195707889299884
    "J8" = 0.3195707889299884
    "M8" = 0.1809866666666667;   "N8" = 0.54296;               "O8" = 0.03987407676082905
    "P8" = 0.03987407676082905;  "Q8" = 2.272096417751111;    "R8" = 20.44886775976
    "S8" = 0.01274259016831306;  "T8" = 0.01274259016831306

    "G9" = 12.55394366666667;    "H9" = 37.661831;             "I9" = 0.3195707889299884
    "J9" = 0.3195707889299884
    "O9" = 0.1057193993302571;   "P9" = 0.1057193993302571;   "Q9" = 6.024081007464112
    "R9" = 54.216729067177;      "S9" = 0.03378483184917475;  "T9" = 0.03378483184917475

    "G10" = 12.55394366666667;   "H10" = 37.661831;            "I10" = 0.3195707889299884
    "J10" = 0.3195707889299884
    "M10" = 3.878113333333333;   "N10" = 11.63434;              "O10" = 0.8544065239089139
    "P10" = 0.8544065239089139;  "Q10" = 48.68561631961555;   "R10" = 438.17054687654
    "S10" = 0.2730433669125006;  "T10" = 0.2730433669125006
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
